$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (Coin name / Link URL) ---
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

# --- Numeric-looking text cells (Price / Volume) ---
# These must stay as text. Excel auto-detects numeric-looking strings
# and converts them to numbers, so for each affected cell we temporarily
# force a Text number format, assign the value, then restore the cell`s
# original style so no visible formatting change is left behind.
# NOTE: this is done per-cell (not via a multi-area Union range) because
# the COM Range setter only actually applies to the first area of a
# multi-area range here.
$numericRiskCells = @(
    'D2',
    'E2',
    'D3',
    'E3',
    'E4',
    'D5',
    'E5',
    'E6',
    'D7',
    'E7',
    'D8',
    'E8',
    'E9',
    'E10',
    'D11',
    'E11',
    'D12',
    'E12',
    'E13',
    'D14',
    'E14',
    'D15',
    'E15',
    'D16',
    'E16',
    'D17',
    'E17',
    'D18',
    'E18',
    'D19',
    'E19',
    'D20',
    'E20',
    'D21',
    'E21',
    'D22',
    'E22',
    'D23',
    'E23',
    'D24',
    'E24',
    'D25',
    'E25',
    'D26',
    'E26',
    'D27',
    'E27',
    'D28',
    'E28',
    'D29',
    'E29',
    'E30',
    'D31',
    'E31',
    'D32',
    'E32',
    'D33',
    'E33',
    'D34',
    'E34',
    'D35',
    'E35',
    'D36',
    'E36',
    'D37',
    'E37',
    'D38',
    'E38',
    'D39',
    'E39',
    'D40',
    'E40',
    'D41',
    'E41',
    'D42',
    'E42',
    'D43',
    'E43',
    'D44',
    'E44',
    'D45',
    'E45',
    'E46',
    'D47',
    'E47',
    'D48',
    'E48',
    'D49',
    'E49',
    'D50',
    'E50',
    'D51',
    'E51'
)

$origStyles = @{}
foreach ($cellRef in $numericRiskCells) {
    $origStyles[$cellRef] = $ws.Range($cellRef).Style
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '96.517.58'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '3.668.17'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '240.16'
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('E6').Value = '  +8.03%  '
$ws.Range('D7').Value = '658.73'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '0.421'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('D11').Value = '3.666.58'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '45.79'
$ws.Range('E12').Value = '  +3.84%  '
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').Value = '6.76'
$ws.Range('E14').Value = '  +4.53%  '
$ws.Range('D15').Value = '4.354.28'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').Value = '0.0000267'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = '96.359.30'
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').Value = '8.87'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('D19').Value = '3.677.38'
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').Value = '12.82'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').Value = '18.61'
$ws.Range('E21').Value = '  +1.88%  '
$ws.Range('D22').Value = '0.525'
$ws.Range('E22').Value = '  -2.61%  '
$ws.Range('D23').Value = '524.24'
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('D24').Value = '3.44'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').Value = '7.07'
$ws.Range('E25').Value = '  +2.30%  '
$ws.Range('D26').Value = '0.0000203'
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('D27').Value = '101.93'
$ws.Range('E27').Value = '  +2.89%  '
$ws.Range('D28').Value = '13.13'
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').Value = '3.870.08'
$ws.Range('E29').Value = '  +1.31%  '
$ws.Range('E30').Value = '  +9.43%  '
$ws.Range('D31').Value = '12.46'
$ws.Range('E31').Value = '  +5.67%  '
$ws.Range('D32').Value = '3.03'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  +14.53%  '
$ws.Range('D35').Value = '0.186'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('D37').Value = '32.48'
$ws.Range('E37').Value = '  +1.93%  '
$ws.Range('D38').Value = '658.13'
$ws.Range('E38').Value = '  +6.14%  '
$ws.Range('D39').Value = '0.589'
$ws.Range('E39').Value = '  +2.67%  '
$ws.Range('D40').Value = '8.88'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').Value = '42.33'
$ws.Range('E41').Value = '  +26.91%  '
$ws.Range('D42').Value = '0.160'
$ws.Range('E42').Value = '  +4.31%  '
$ws.Range('D43').Value = '1.98'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('D44').Value = '0.961'
$ws.Range('E44').Value = '  +2.86%  '
$ws.Range('D45').Value = '6.51'
$ws.Range('E45').Value = '  +8.89%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').Value = '0.0467'
$ws.Range('E47').Value = '  +5.94%  '
$ws.Range('D48').Value = '0.449'
$ws.Range('E48').Value = '  +15.71%  '
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').Value = '23.64'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').Value = '8.63'
$ws.Range('E51').Value = '  +0.67%  '

foreach ($cellRef in $numericRiskCells) {
    $ws.Range($cellRef).Style = $origStyles[$cellRef]
}
